# Update column F (dSF) values on Sheet1 to reflect the re-pulled /
# recalculated data ("repull data, push all data, mean calculation").
# Only the dSF column (column F) changes; all other columns/rows are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = -4
$ws.Cells.Item(3, 6).Value = -1
$ws.Cells.Item(4, 6).Value = -6
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(8, 6).Value = -2
$ws.Cells.Item(9, 6).Value = -1
$ws.Cells.Item(10, 6).Value = 4
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(13, 6).Value = -1
$ws.Cells.Item(14, 6).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(16, 6).Value = -1
$ws.Cells.Item(17, 6).Value = 12
$ws.Cells.Item(18, 6).Value = -7
$ws.Cells.Item(19, 6).Value = -6
$ws.Cells.Item(21, 6).Value = -1
$ws.Cells.Item(23, 6).Value = 5
$ws.Cells.Item(25, 6).Value = -2
$ws.Cells.Item(26, 6).Value = -1
$ws.Cells.Item(27, 6).Value = -4
$ws.Cells.Item(28, 6).Value = -3
$ws.Cells.Item(29, 6).Value = -1
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(31, 6).Value = 2
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(34, 6).Value = 2
$ws.Cells.Item(35, 6).Value = -6
$ws.Cells.Item(36, 6).Value = -6
$ws.Cells.Item(39, 6).Value = -2
$ws.Cells.Item(40, 6).Value = -1
$ws.Cells.Item(41, 6).Value = -2
